# Insert a new weekly price-record row for "Ajo" (Chino / Primera, $/malla 10 kilos)
# at row 73, shifting the existing rows 73:174 down to 74:175.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("73:73").Insert()

$ws.Cells.Item(73, 1).Value = 9
$ws.Cells.Item(73, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(73, 3).Value = "Metropolitana"
$ws.Cells.Item(73, 4).Value = 44571
$ws.Cells.Item(73, 5).Value = 13
$ws.Cells.Item(73, 6).Value = 100112003
$ws.Cells.Item(73, 7).Value = "Ajo"
$ws.Cells.Item(73, 8).Value = "Chino"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 340
$ws.Cells.Item(73, 11).Value = 18000
$ws.Cells.Item(73, 12).Value = 18500
$ws.Cells.Item(73, 13).Value = 18250
$ws.Cells.Item(73, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(73, 15).Value = "China"
$ws.Cells.Item(73, 16).Value = 1825
$ws.Cells.Item(73, 17).Value = 10
$ws.Cells.Item(73, 18).Value = "Hortaliza"
